$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.223.00'
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("D3").Value = '1.896.14'
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").Value = "'245.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").Value = "'0.683"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.17%  '
$ws.Range("E7").Value = '  -0.35%  '
$ws.Range("D8").Value = "'40.79"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.36%  '
$ws.Range("D9").Value = "'0.346"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.69%  '
$ws.Range("D10").Value = "'52.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.72%  '
$ws.Range("E11").Value = '  +1.99%  '
$ws.Range("D12").Value = "'0.0982"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.54%  '
$ws.Range("D13").Value = '2.170.37'
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("D14").Value = "'12.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.91%  '
$ws.Range("D15").Value = "'0.702"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.57%  '
$ws.Range("D16").Value = '1.895.87'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("D18").Value = '35.178.09'
$ws.Range("E18").Value = '  -1.04%  '
$ws.Range("D19").Value = "'71.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").Value = '0.0₃0815'
$ws.Range("E20").Value = '  +0.70%  '
$ws.Range("D21").Value = "'240.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.40%  '
$ws.Range("D22").Value = "'12.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.84%  '
$ws.Range("D23").Value = "'4.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.42%  '
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").Value = "'2.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.54%  '
$ws.Range("E26").Value = '  +3.49%  '
$ws.Range("D27").Value = "'167.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.52%  '
$ws.Range("D28").Value = "'8.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").Value = "'18.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.90%  '
$ws.Range("E30").Value = '  +3.39%  '
$ws.Range("E32").Value = '  +1.12%  '
$ws.Range("D33").Value = "'0.0566"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("D34").Value = "'1.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.42%  '
$ws.Range("E35").Value = '  -0.44%  '
$ws.Range("E36").Value = '  -7.76%  '
$ws.Range("D37").Value = "'4.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("D38").Value = "'1.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.99%  '
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("D40").Value = "'16.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.99%  '
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("E42").Value = '  +0.65%  '
$ws.Range("D43").Value = "'0.0635"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.38%  '
$ws.Range("D44").Value = "'89.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.81%  '
$ws.Range("D45").Value = '1.348.02'
$ws.Range("E45").Value = '  -0.43%  '
$ws.Range("E46").Value = '  +2.19%  '
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").Value = "'2.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.95%  '
$ws.Range("D49").Value = "'45.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -14.05%  '
$ws.Range("E50").Value = '  -5.07%  '
$ws.Range("D51").Value = "'6.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.04%  '
